$d = $word.ActiveDocument

function Wrap-PkgXml($innerBodyXml) {
    return '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Step 1: split "Register del producer" -> "Register" (proofErr) + " del producer"
# within the first list paragraph (keeping the other runs of that paragraph intact).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(2)
$body = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Register</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> del producer</w:t></w:r>' +
    '<w:r><w:t>, quindi s</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">i va a popolare una riga della tabella </w:t></w:r>' +
    '<w:r><w:t>db producer</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> andando a impostare:</w:t></w:r>' +
    '</w:p></w:body>'
$p.Range.InsertXML((Wrap-PkgXml $body))

# ---------------------------------------------------------------------------
# Step 2: split "tetti massimi x ogni slot" -> "tetti massimi" + " " + "x" (proofErr) + " ogni slot"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(3)
$body = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>tetti massimi</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>x</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> ogni slot</w:t></w:r>' +
    '</w:p></w:body>'
$p.Range.InsertXML((Wrap-PkgXml $body))

# ---------------------------------------------------------------------------
# Step 3: split "register del consumer" -> "register" (proofErr) + " del consumer"
# (this paragraph does NOT get strikethrough)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(9)
$body = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>register</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> del consumer</w:t></w:r>' +
    '</w:p></w:body>'
$p.Range.InsertXML((Wrap-PkgXml $body))

# ---------------------------------------------------------------------------
# Step 4: apply strikethrough formatting to paragraphs 2-8
# (Register del producer... through login del producer...)
# ---------------------------------------------------------------------------
for ($i = 2; $i -le 8; $i++) {
    $d.Paragraphs($i).Range.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------------
# Step 5: insert two new list paragraphs ("register admin", "login admin")
# right after "login del consumer" (paragraph 10) and before
# "il consumer opziona..." (paragraph 11).
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$p10.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(11)
$body = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>register</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> admin</w:t></w:r>' +
    '</w:p></w:body>'
$newPara.Range.InsertXML((Wrap-PkgXml $body))

$p11 = $d.Paragraphs(11)
$p11.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs(12)
$body = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>login admin</w:t></w:r>' +
    '</w:p></w:body>'
$newPara2.Range.InsertXML((Wrap-PkgXml $body))

# ---------------------------------------------------------------------------
# Step 6: insert a new list paragraph ("admin fa la ricarica al consumer")
# right after "si genera una riga in db storico acquisti settando i vari campi"
# and before the trailing empty paragraph.
# ---------------------------------------------------------------------------
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($ptext -eq "si genera una riga in db storico acquisti settando i vari campi") {
        $targetIdx = $i
    }
}
$pTarget = $d.Paragraphs($targetIdx)
$pTarget.Range.InsertParagraphAfter()
$newPara3 = $d.Paragraphs($targetIdx + 1)
$body = '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>admin fa la ricarica al consumer</w:t></w:r>' +
    '</w:p></w:body>'
$newPara3.Range.InsertXML((Wrap-PkgXml $body))

Write-Output "All edits applied"
